$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "password" column (column E) -- shifts games/hosted games/role left
$ws.Columns("E").Delete()

# After the shift: E2 holds the old "games" value ("0;"), F2 holds the old
# "hosted games" value ("0;1;2;"). Move the games value into the hosted games
# cell (overwriting it) and clear the games cell, matching the target state.
$ws.Range("E2").Cut($ws.Range("F2"))
